# Update to LDkit 2.0.0 results: new raw measurements for the LDkit
# "Search by field" scenario (columns K and L on the Measurements sheet),
# a relabelled header (x100 -> x1000) and the corresponding normalisation
# divisor on the "Normalised measurements" sheet.

$wb = $excel.ActiveWorkbook
$wsMeasurements = $wb.Worksheets.Item("Measurements")
$wsNormalised   = $wb.Worksheets.Item("Normalised measurements")

# New raw values for columns K (LDkit - Search by field, x100) and
# L (LDkit - Search by field, now x1000) for rows 3..32.
$newData = @(
    @(3,  15703, 6200),
    @(4,  14873, 5741),
    @(5,  16505, 6383),
    @(6,  14648, 6758),
    @(7,  15209, 6028),
    @(8,  14891, 5029),
    @(9,  16938, 5233),
    @(10, 20885, 5632),
    @(11, 19795, 6707),
    @(12, 17525, 6791),
    @(13, 16201, 6839),
    @(14, 14705, 6603),
    @(15, 14603, 6543),
    @(16, 14719, 6527),
    @(17, 14712, 7027),
    @(18, 14945, 6739),
    @(19, 15761, 6951),
    @(20, 15864, 6721),
    @(21, 15921, 6996),
    @(22, 15613, 6906),
    @(23, 15928, 6890),
    @(24, 15741, 6569),
    @(25, 15753, 6504),
    @(26, 15855, 6895),
    @(27, 16114, 6809),
    @(28, 16246, 6646),
    @(29, 16188, 6638),
    @(30, 16268, 6863),
    @(31, 15639, 6577),
    @(32, 15568, 7061)
)

foreach ($entry in $newData) {
    $row = $entry[0]
    $kVal = $entry[1]
    $lVal = $entry[2]
    $wsMeasurements.Cells.Item($row, 11).Value = $kVal
    $wsMeasurements.Cells.Item($row, 12).Value = $lVal
}

# Header relabel: LDkit's "Search by field" column is now reported as
# x1000 instead of x100 (both sheets share the same header layout).
$wsMeasurements.Range("L2").Value = "Search by field (x1000)"
$wsNormalised.Range("L2").Value = "Search by field (x1000)"

# The normalisation formulas for column L (LDkit "Search by field")
# now divide by 1000 instead of 100 to match the new column header.
for ($row = 3; $row -le 32; $row++) {
    $wsNormalised.Cells.Item($row, 12).Formula = "=Measurements!L$row/1000"
}
